$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.991.94"
$ws.Range("E2").Value = "  +1.34%  "

# Row 3
$ws.Range("D3").Value = "1.675.35"
$ws.Range("E3").Value = "  +2.47%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "215.95"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6
$ws.Range("D6").Value = "0.530"
$ws.Range("E6").Value = "  +5.39%  "

# Row 7
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +3.08%  "

# Row 9
$ws.Range("D9").Value = "0.0621"
$ws.Range("E9").Value = "  +1.57%  "

# Row 10
$ws.Range("D10").Value = "20.37"
$ws.Range("E10").Value = "  +5.34%  "

# Row 11
$ws.Range("E11").Value = "  +3.92%  "

# Row 12
$ws.Range("D12").Value = "1.902.75"
$ws.Range("E12").Value = "  +2.21%  "

# Row 13
$ws.Range("D13").Value = "1.678.55"
$ws.Range("E13").Value = "  +2.72%  "

# Row 14
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  +0.93%  "

# Row 15
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  +1.72%  "

# Row 16
$ws.Range("D16").Value = "65.78"
$ws.Range("E16").Value = "  +2.64%  "

# Row 17
$ws.Range("D17").Value = "26.995.08"
$ws.Range("E17").Value = "  +1.35%  "

# Row 18
$ws.Range("D18").Value = "233.15"
$ws.Range("E18").Value = "  -1.05%  "

# Row 19
$ws.Range("D19").Value = "7.83"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0736"
$ws.Range("E20").Value = "  +1.44%  "

# Row 21
$ws.Range("E21").Value = "  +0.15%  "

# Row 22
$ws.Range("D22").Value = "4.46"
$ws.Range("E22").Value = "  +2.77%  "

# Row 23
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "9.24"
$ws.Range("E23").Value = "  +0.54%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").Value = "  +0.32%  "

# Row 25
$ws.Range("D25").Value = "145.56"
$ws.Range("E25").Value = "  -0.53%  "

# Row 26
$ws.Range("E26").Value = "  +0.78%  "

# Row 27
$ws.Range("E27").Value = "  +2.57%  "

# Row 28
$ws.Range("D28").Value = "16.04"
$ws.Range("E28").Value = "  +1.89%  "

# Row 29
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("E30").Value = "  +0.12%  "

# Row 31
$ws.Range("E31").Value = "  +1.51%  "

# Row 32
$ws.Range("E32").Value = "  +1.61%  "

# Row 33
$ws.Range("D33").Value = "1.458.21"
$ws.Range("E33").Value = "  -4.44%  "

# Row 34
$ws.Range("E34").Value = "  +4.56%  "

# Row 35
$ws.Range("D35").Value = "1.62"
$ws.Range("E35").Value = "  +5.53%  "

# Row 36
$ws.Range("E36").Value = "  -0.02%  "

# Row 37
$ws.Range("D37").Value = "0.908"
$ws.Range("E37").Value = "  +7.70%  "

# Row 38
$ws.Range("D38").Value = "0.569"
$ws.Range("E38").Value = "  -0.51%  "

# Row 39
$ws.Range("E39").Value = "  +1.03%  "

# Row 40
$ws.Range("D40").Value = "6.03"
$ws.Range("E40").Value = "  +2.14%  "

# Row 41
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.07%  "

# Row 42
$ws.Range("D42").Value = "2.30"
$ws.Range("E42").Value = "  +3.68%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "65.94"
$ws.Range("E43").Value = "  +3.64%  "

# Row 44
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "0.970"
$ws.Range("E44").Value = "  +6.76%  "

# Row 45
$ws.Range("D45").Value = "1.812.09"
$ws.Range("E45").Value = "  +2.17%  "

# Row 46
$ws.Range("D46").Value = "0.782"
$ws.Range("E46").Value = "  +2.50%  "

# Row 47
$ws.Range("D47").Value = "90.60"
$ws.Range("E47").Value = "  +0.21%  "

# Row 48
$ws.Range("D48").Value = "1.54"
$ws.Range("E48").Value = "  +0.96%  "

# Row 49
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  +3.32%  "

# Row 50
$ws.Range("D50").Value = "0.0508"
$ws.Range("E50").Value = "  +1.31%  "

# Row 51
$ws.Range("D51").Value = "7.61"
$ws.Range("E51").Value = "  +0.11%  "
